$wb = $excel.ActiveWorkbook

# The active/selected sheet is "VEDA_Sets-Proc" (tab index 2, 1-based)
$ws = $wb.Worksheets.Item("VEDA_Sets-Proc")

# Append a new data row (row 21) with the new PSET entry
$ws.Range("F21").Value = "Grid"
$ws.Range("A21").Value = "IRE"
$ws.Range("B21").Value = "g[_]*"

# Reflect the new active cell/selection on that sheet
$ws.Activate()
$ws.Range("B21").Select()
